{"js": "// Replace three-digit x one-digit multiplication equations with new values.\n// Each old equation text is unique in the document, so a direct search+replace is safe.\nconst replacements = [\n  [\"637\u00d72=1274\", \"523\u00d75=2615\"],\n  [\"947\u00d73=2841\", \"316\u00d78=2528\"],\n  [\"266\u00d73=798\", \"534\u00d79=4806\"],\n  [\"203\u00d73=609\", \"800\u00d78=6400\"],\n  [\"710\u00d76=4260\", \"554\u00d77=3878\"],\n  [\"375\u00d77=2625\", \"518\u00d72=1036\"],\n  [\"159\u00d77=1113\", \"875\u00d75=4375\"],\n  [\"145\u00d75=725\", \"571\u00d75=2855\"],\n  [\"928\u00d79=8352\", \"951\u00d76=5706\"],\n  [\"891\u00d78=7128\", \"478\u00d75=2390\"],\n  [\"141\u00d75=705\", \"685\u00d75=3425\"],\n  [\"510\u00d74=2040\", \"442\u00d77=3094\"],\n  [\"546\u00d75=2730\", \"887\u00d74=3548\"],\n  [\"683\u00d79=6147\", \"387\u00d72=774\"],\n  [\"521\u00d72=1042\", \"206\u00d75=1030\"],\n  [\"566\u00d74=2264\", \"522\u00d74=2088\"],\n  [\"186\u00d74=744\", \"442\u00d74=1768\"],\n  [\"439\u00d77=3073\", \"525\u00d76=3150\"],\n  [\"254\u00d75=1270\", \"726\u00d77=5082\"],\n  [\"989\u00d73=2967\", \"601\u00d73=1803\"],\n  [\"848\u00d79=7632\", \"219\u00d75=1095\"],\n  [\"649\u00d79=5841\", \"297\u00d77=2079\"],\n  [\"705\u00d79=6345\", \"414\u00d78=3312\"],\n  [\"144\u00d79=1296\", \"102\u00d73=306\"],\n  [\"918\u00d77=6426\", \"898\u00d75=4490\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "# Replace three-digit x one-digit multiplication equations with new values.\n# Each old equation text occurs exactly once in the document, so Find/Replace\n# per pair (without ReplaceAll ambiguity) is safe and deterministic.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"637\u00d72=1274\"; New = \"523\u00d75=2615\"},\n    @{Old = \"947\u00d73=2841\"; New = \"316\u00d78=2528\"},\n    @{Old = \"266\u00d73=798\"; New = \"534\u00d79=4806\"},\n    @{Old = \"203\u00d73=609\"; New = \"800\u00d78=6400\"},\n    @{Old = \"710\u00d76=4260\"; New = \"554\u00d77=3878\"},\n    @{Old = \"375\u00d77=2625\"; New = \"518\u00d72=1036\"},\n    @{Old = \"159\u00d77=1113\"; New = \"875\u00d75=4375\"},\n    @{Old = \"145\u00d75=725\"; New = \"571\u00d75=2855\"},\n    @{Old = \"928\u00d79=8352\"; New = \"951\u00d76=5706\"},\n    @{Old = \"891\u00d78=7128\"; New = \"478\u00d75=2390\"},\n    @{Old = \"141\u00d75=705\"; New = \"685\u00d75=3425\"},\n    @{Old = \"510\u00d74=2040\"; New = \"442\u00d77=3094\"},\n    @{Old = \"546\u00d75=2730\"; New = \"887\u00d74=3548\"},\n    @{Old = \"683\u00d79=6147\"; New = \"387\u00d72=774\"},\n    @{Old = \"521\u00d72=1042\"; New = \"206\u00d75=1030\"},\n    @{Old = \"566\u00d74=2264\"; New = \"522\u00d74=2088\"},\n    @{Old = \"186\u00d74=744\"; New = \"442\u00d74=1768\"},\n    @{Old = \"439\u00d77=3073\"; New = \"525\u00d76=3150\"},\n    @{Old = \"254\u00d75=1270\"; New = \"726\u00d77=5082\"},\n    @{Old = \"989\u00d73=2967\"; New = \"601\u00d73=1803\"},\n    @{Old = \"848\u00d79=7632\"; New = \"219\u00d75=1095\"},\n    @{Old = \"649\u00d79=5841\"; New = \"297\u00d77=2079\"},\n    @{Old = \"705\u00d79=6345\"; New = \"414\u00d78=3312\"},\n    @{Old = \"144\u00d79=1296\"; New = \"102\u00d73=306\"},\n    @{Old = \"918\u00d77=6426\"; New = \"898\u00d75=4490\"}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}"}
